# edit.ps1 - applies the "streaming-figures" deck update:
#   1) Refresh the cached "datetimeFigureOut" date field text (1/23/14 -> 9/6/14)
#      on the slide master and every slide layout's Date Placeholder.
#   2) On slide 1's architecture diagram, relabel the "HDFS" source box to
#      "HDFS/S3" and the "ZeroMQ" source box to "Kinesis".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholders: slide master + all custom (slide) layouts.
# ---------------------------------------------------------------------------
$newDate = "9/6/14"

$sm = $p.SlideMaster
for ($i = 1; $i -le $sm.Shapes.Count; $i++) {
    $sh = $sm.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($L = 1; $L -le $sm.CustomLayouts.Count; $L++) {
    $layout = $sm.CustomLayouts.Item($L)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 1 streaming-sources diagram: HDFS -> HDFS/S3, ZeroMQ -> Kinesis.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$rootGroup = $slide1.Shapes.Item(1)

$hdfsShape = $rootGroup.GroupItems.Item("Rounded Rectangle 47")
$hdfsShape.TextFrame.TextRange.Text = "HDFS/S3"

$zeroMqShape = $rootGroup.GroupItems.Item("Rounded Rectangle 48")
$zeroMqShape.TextFrame.TextRange.Text = "Kinesis"
